# GPLIM-2957: fix excel headers so they match with values in Decision
#
# The manifest's header row used column titles that no longer line up
# with the values the Decision import expects:
#   A1: "Sample ID"  -> "Specimen_Number"
#   F1: "T/N"         -> "SAMPLE_TYPE"
#
# The leading "'" forces Excel to keep treating the header cells as
# quoted/text entries (matching their original quote-prefixed cell
# style) instead of re-evaluating them as something that could be
# parsed differently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "SAMPLE_TYPE"
$ws.Range("A1").Value = "'Specimen_Number"
